# Loan RBI, Variable Instalments
# Insert a new (blank) column into the "Repayment Schedule" sheet before
# column N ("Late"), shifting the "Late" / "heading" / "Outstanding"
# columns one place to the right, then leave the sheet active with the
# selection where the user last clicked.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a blank column before column N - everything from N onward shifts
# right by one (N->O, O->P, P->Q), matching formatting/styles of the
# column being pushed aside.
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab and select the cell the
# user ended up on after the edit.
$ws.Activate()
$ws.Range("K16").Select()
